# Applies the "added missing cal events, assigned OOI bar codes where
# necessary, corrected integration events" edit to the Asset_Cal_Info sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# Row 11 (the GL494 "-00-ENG000000" integration event): the Sensor OOIBARCODE
# (col E) was missing and the Sensor Serial Number (col F) incorrectly held
# the bare deployment number (494) instead of the actual controller asset.
$ws.Range("E11").Value = "OL000109"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = "GL494 Controller"
